# Added systematic uncertainty to ball drop
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terrence")

# Make 4 blank rows of room before the old "Height of comb" block (old row 27),
# shifting everything below down by four rows (27->31, 29->33, 31->35, ... 67->71).
$ws.Rows("26:29").Insert()

# --- New "systematic uncertainty" columns (D/I/N/S/X) for each of the five runs ---
# Row 6 (+0.5% height / +0.5% error uncertainty variant)
$ws.Range("D6").Formula  = "=(2/(A6+B7))*((0.4545/B7)-(0.0795/A6))"
$ws.Range("I6").Formula  = "=(2/(F6 + G7)) * ((0.5395/G7) - (0.0795/F6))"
$ws.Range("N6").Formula  = "=(2/(K6+L7))*((0.6295/L7)-(0.0795/K6))"
$ws.Range("S6").Formula  = "=(2/(P6+Q7))*((0.2695/Q7)-(0.0795/P6))"
$ws.Range("X6").Formula  = "=(2/(U6+V7))*((0.7395/V7)-(0.0795/U6))"

# Row 7 (-0.5% height / -0.5% error uncertainty variant)
$ws.Range("D7").Formula  = "=(2/(A6 + B7)) * ((0.4555/B7) - (0.0805/A6))"
$ws.Range("I7").Formula  = "=(2/(F6 + G7)) * ((0.5405/G7) - (0.0805/F6))"
$ws.Range("N7").Formula  = "=(2/(K6 + L7)) * ((0.6305/L7) - (0.0805/K6))"
$ws.Range("S7").Formula  = "=(2/(P6 + Q7)) * ((0.2705/Q7) - (0.0805/P6))"
$ws.Range("X7").Formula  = "=(2/(U6 + V7)) * ((0.7405/V7) - (0.0805/U6))"

# --- Row 26: Average ---
$ws.Range("C26").Value = "Average"
$ws.Range("D26").Formula = "=AVERAGE(C6,C8,C10,C12,C14,C16,C18,C20,C22,C24)"
$ws.Range("H26").Value = "Average"
$ws.Range("I26").Formula = "=AVERAGE(H6,H8,H10,H12,H14,H16,H18,H20,H22,H24)"
$ws.Range("M26").Value = "Average"
$ws.Range("N26").Formula = "=AVERAGE(M6,M8,M10,M12,M14,M16,M18,M20,M22,M24)"
$ws.Range("R26").Value = "Average"
$ws.Range("S26").Formula = "=AVERAGE(R6,R8,R10,R12,R14,R16,R18,R20,R22,R24)"
$ws.Range("W26").Value = "Average"
$ws.Range("X26").Formula = "=AVERAGE(W6,W8,W10,W12,W14,W16,W18,W20,W22,W24)"

# --- Row 27: Statistical Uncertainty ---
$ws.Range("C27").Value = "Statistical Uncertainty"
$ws.Range("D27").Formula = "=STDEV.S(C6,C8,C10,C12,C14,C16,C18,C20,C22,C24,D26)/SQRT(10)"
$ws.Range("H27").Value = "Statistical Uncertainty"
$ws.Range("I27").Formula = "=STDEV.S(H6,H8,H10,H12,H14,H16,H18,H20,H22,H24)/SQRT(10)"
$ws.Range("M27").Value = "Statistical Uncertainty"
$ws.Range("N27").Formula = "=STDEV.S(M6,M8,M10,M12,M14,M16,M18,M20,M22,M24)/SQRT(10)"
$ws.Range("R27").Value = "Statistical Uncertainty"
$ws.Range("S27").Formula = "=STDEV.S(R6,R8,R10,R12,R14,R16,R18,R20,R22,R24)/SQRT(10)"
$ws.Range("W27").Value = "Statistical Uncertainty"
$ws.Range("X27").Formula = "=STDEV.S(W6,W8,W10,W12,W14,W16,W18,W20,W22,W24)/SQRT(10)"

# --- Row 28: Systematic Uncertainty ---
$ws.Range("C28").Value = "Systematic Uncertainty"
$ws.Range("D28").Formula = "=(D6-D7)/2"
$ws.Range("H28").Value = "Systematic Uncertainty"
$ws.Range("I28").Formula = "=(I6 - I7)/2"
$ws.Range("M28").Value = "Systematic Uncertainty"
$ws.Range("N28").Formula = "=(N6-N7)/2"
$ws.Range("R28").Value = "Systematic Uncertainty"
$ws.Range("S28").Formula = "=(S6-S7)/2"
$ws.Range("W28").Value = "Systematic Uncertainty"
$ws.Range("X28").Formula = "=(X6-X7)/2"

# --- Restore the view the author left the sheet in ---
$ws.Range("D7").Select()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("N29").Select()
